$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting text for the R10 rule row
$ws.Range("E8").Value = "GIT UPDATE"

# Move the active selection to E8 (matches the saved selection in the file)
$ws.Activate()
$ws.Range("E8").Select()
